$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B91 currently holds the text "3" as a string; convert it to a real number
$ws.Range("B91").Value = 3

# Append new row 92 with the new annotation data
$ws.Range("A92").Value = "Ruilin"
$ws.Range("B92").Value = "'2"
$ws.Range("B92").Style = "Normal"
$ws.Range("C92").Value = "Authors, please post a rebuttal soon if you are planning on it."
$ws.Range("D92").Value = "CRT"
$ws.Range("E92").Value = "OTH"
$ws.Range("F92").Value = "0e22fe07-2d2d-417e-8066-2728b416bb18"
$ws.Range("G92").Value = "Byht0GbRZ_annotated.xlsx"
$ws.Range("H92").Value = "Authors, please post a rebuttal soon if you are planning on it."
